# Auto-generated edit script: updates cached market-price values (columns H-N)
# on several leve-profit worksheets, per scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 11000
$ws.Cells.Item(13, 10).Value = 11000
$ws.Cells.Item(13, 12).Value = 11000
$ws.Cells.Item(13, 14).Value = -11338
$ws.Cells.Item(15, 8).Value = 1150.4861
$ws.Cells.Item(15, 9).Value = 1150.4861
$ws.Cells.Item(15, 11).Value = 3451.4583
$ws.Cells.Item(15, 13).Value = -3282.4583
$ws.Cells.Item(53, 8).Value = 380.2903
$ws.Cells.Item(53, 9).Value = 402.66666
$ws.Cells.Item(53, 10).Value = 366.1579
$ws.Cells.Item(53, 11).Value = 402.66666
$ws.Cells.Item(53, 12).Value = 366.1579
$ws.Cells.Item(53, 13).Value = 234.33334
$ws.Cells.Item(53, 14).Value = -1640.1579
$ws.Cells.Item(106, 8).Value = 2267.2727
$ws.Cells.Item(106, 9).Value = 2117.5
$ws.Cells.Item(106, 10).Value = 2666.6667
$ws.Cells.Item(106, 11).Value = 2117.5
$ws.Cells.Item(106, 12).Value = 2666.6667
$ws.Cells.Item(106, 13).Value = -1486.5
$ws.Cells.Item(106, 14).Value = -3928.6667
$ws.Cells.Item(129, 8).Value = 3131.2222
$ws.Cells.Item(129, 10).Value = 1061.4062
$ws.Cells.Item(129, 12).Value = 3184.2186
$ws.Cells.Item(129, 14).Value = -13184.2186
$ws.Cells.Item(137, 8).Value = 1809.28
$ws.Cells.Item(137, 9).Value = 1549.1052
$ws.Cells.Item(137, 10).Value = 2633.1667
$ws.Cells.Item(137, 11).Value = 4647.3156
$ws.Cells.Item(137, 12).Value = 7899.500100000001
$ws.Cells.Item(137, 13).Value = -2097.3156
$ws.Cells.Item(137, 14).Value = -12999.5001
$ws.Cells.Item(138, 8).Value = 2129.7666
$ws.Cells.Item(138, 9).Value = 1830.5714
$ws.Cells.Item(138, 10).Value = 2827.889
$ws.Cells.Item(138, 11).Value = 5491.7142
$ws.Cells.Item(138, 12).Value = 8483.667000000001
$ws.Cells.Item(138, 13).Value = -351.7142000000003
$ws.Cells.Item(138, 14).Value = -18763.667
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 2073
$ws.Cells.Item(26, 9).Value = 2073
$ws.Cells.Item(26, 10).Value = 0
$ws.Cells.Item(26, 11).Value = 2073
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 13).Value = -1743
$ws.Cells.Item(26, 14).ClearContents()
$ws.Cells.Item(32, 8).Value = 20898.822
$ws.Cells.Item(32, 9).Value = 4618.101
$ws.Cells.Item(32, 10).Value = 137824
$ws.Cells.Item(32, 11).Value = 4618.101
$ws.Cells.Item(32, 12).Value = 137824
$ws.Cells.Item(32, 13).Value = -4331.101
$ws.Cells.Item(32, 14).Value = -138398
$ws.Cells.Item(39, 8).Value = 5100
$ws.Cells.Item(39, 9).Value = 5100
$ws.Cells.Item(39, 11).Value = 5100
$ws.Cells.Item(39, 13).Value = -4580
$ws.Cells.Item(61, 8).Value = 2206.0698
$ws.Cells.Item(61, 9).Value = 1291.3158
$ws.Cells.Item(61, 10).Value = 2930.25
$ws.Cells.Item(61, 11).Value = 1291.3158
$ws.Cells.Item(61, 12).Value = 2930.25
$ws.Cells.Item(61, 13).Value = -1079.3158
$ws.Cells.Item(61, 14).Value = -3354.25
$ws.Cells.Item(74, 8).Value = 764.63635
$ws.Cells.Item(74, 9).Value = 756.7778
$ws.Cells.Item(74, 10).Value = 800
$ws.Cells.Item(74, 11).Value = 756.7778
$ws.Cells.Item(74, 12).Value = 800
$ws.Cells.Item(74, 13).Value = 117.2222
$ws.Cells.Item(74, 14).Value = -2548
$ws.Cells.Item(77, 8).Value = 764.63635
$ws.Cells.Item(77, 9).Value = 756.7778
$ws.Cells.Item(77, 10).Value = 800
$ws.Cells.Item(77, 11).Value = 3783.889
$ws.Cells.Item(77, 12).Value = 4000
$ws.Cells.Item(77, 13).Value = 584.1110000000003
$ws.Cells.Item(77, 14).Value = -12736
$ws.Cells.Item(132, 8).Value = 2962.6128
$ws.Cells.Item(132, 9).Value = 2569.682
$ws.Cells.Item(132, 10).Value = 3923.111
$ws.Cells.Item(132, 11).Value = 7709.045999999999
$ws.Cells.Item(132, 12).Value = 11769.333
$ws.Cells.Item(132, 13).Value = -5179.045999999999
$ws.Cells.Item(132, 14).Value = -16829.333
$ws.Cells.Item(136, 8).Value = 2206.0698
$ws.Cells.Item(136, 9).Value = 1291.3158
$ws.Cells.Item(136, 10).Value = 2930.25
$ws.Cells.Item(136, 11).Value = 3873.9474
$ws.Cells.Item(136, 12).Value = 8790.75
$ws.Cells.Item(136, 13).Value = -1323.9474
$ws.Cells.Item(136, 14).Value = -13890.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 14).ClearContents()
$ws.Cells.Item(64, 8).Value = 699
$ws.Cells.Item(64, 9).Value = 201.5
$ws.Cells.Item(64, 10).Value = 983.2857
$ws.Cells.Item(64, 11).Value = 201.5
$ws.Cells.Item(64, 12).Value = 983.2857
$ws.Cells.Item(64, 13).Value = 23.5
$ws.Cells.Item(64, 14).Value = -1433.2857
$ws.Cells.Item(67, 8).Value = 699
$ws.Cells.Item(67, 9).Value = 201.5
$ws.Cells.Item(67, 10).Value = 983.2857
$ws.Cells.Item(67, 11).Value = 201.5
$ws.Cells.Item(67, 12).Value = 983.2857
$ws.Cells.Item(67, 13).Value = 578.5
$ws.Cells.Item(67, 14).Value = -2543.2857
$ws.Cells.Item(94, 8).Value = 497
$ws.Cells.Item(94, 9).Value = 350.2857
$ws.Cells.Item(94, 11).Value = 350.2857
$ws.Cells.Item(94, 13).Value = 100.7143
$ws.Cells.Item(134, 8).Value = 5635.9395
$ws.Cells.Item(134, 9).Value = 5611.826
$ws.Cells.Item(134, 10).Value = 5691.4
$ws.Cells.Item(134, 11).Value = 16835.478
$ws.Cells.Item(134, 12).Value = 17074.2
$ws.Cells.Item(134, 13).Value = -14300.478
$ws.Cells.Item(134, 14).Value = -22144.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 23797.938
$ws.Cells.Item(31, 9).Value = 1094.9697
$ws.Cells.Item(31, 10).Value = 47965.613
$ws.Cells.Item(31, 11).Value = 1094.9697
$ws.Cells.Item(31, 12).Value = 47965.613
$ws.Cells.Item(31, 13).Value = -799.9697000000001
$ws.Cells.Item(31, 14).Value = -48555.613
$ws.Cells.Item(34, 8).Value = 23797.938
$ws.Cells.Item(34, 9).Value = 1094.9697
$ws.Cells.Item(34, 10).Value = 47965.613
$ws.Cells.Item(34, 11).Value = 1094.9697
$ws.Cells.Item(34, 12).Value = 47965.613
$ws.Cells.Item(34, 13).Value = -892.9697000000001
$ws.Cells.Item(34, 14).Value = -48369.613
$ws.Cells.Item(86, 8).Value = 3075.5
$ws.Cells.Item(86, 10).Value = 3067.4443
$ws.Cells.Item(86, 12).Value = 3067.4443
$ws.Cells.Item(86, 14).Value = -5313.4443
$ws.Cells.Item(89, 8).Value = 3075.5
$ws.Cells.Item(89, 10).Value = 3067.4443
$ws.Cells.Item(89, 12).Value = 15337.2215
$ws.Cells.Item(89, 14).Value = -26569.2215
$ws.Cells.Item(107, 8).Value = 868.7895
$ws.Cells.Item(107, 9).Value = 999.5
$ws.Cells.Item(107, 10).Value = 644.7143
$ws.Cells.Item(107, 11).Value = 999.5
$ws.Cells.Item(107, 12).Value = 644.7143
$ws.Cells.Item(107, 13).Value = 920.5
$ws.Cells.Item(107, 14).Value = -4484.7143
$ws.Cells.Item(132, 8).Value = 34887276
$ws.Cells.Item(132, 9).Value = 33336846
$ws.Cells.Item(132, 10).Value = 38465190
$ws.Cells.Item(132, 11).Value = 100010538
$ws.Cells.Item(132, 12).Value = 115395570
$ws.Cells.Item(132, 13).Value = -100008008
$ws.Cells.Item(132, 14).Value = -115400630
$ws.Cells.Item(134, 8).Value = 1150.4878
$ws.Cells.Item(134, 9).Value = 1032.037
$ws.Cells.Item(134, 10).Value = 1378.9286
$ws.Cells.Item(134, 11).Value = 3096.111
$ws.Cells.Item(134, 12).Value = 4136.7858
$ws.Cells.Item(134, 13).Value = -561.1109999999999
$ws.Cells.Item(134, 14).Value = -9206.7858
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 1540.721
$ws.Cells.Item(131, 10).Value = 1522.9156
$ws.Cells.Item(131, 12).Value = 4568.7468
$ws.Cells.Item(131, 14).Value = -14648.7468
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 777353.4
$ws.Cells.Item(107, 10).Value = 10101010
$ws.Cells.Item(107, 12).Value = 10101010
$ws.Cells.Item(107, 14).Value = -10104850
$ws.Cells.Item(113, 8).Value = 1460.5
$ws.Cells.Item(113, 9).Value = 830.5
$ws.Cells.Item(113, 11).Value = 830.5
$ws.Cells.Item(113, 13).Value = 1339.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4048
$ws.Cells.Item(7, 9).Value = 3000
$ws.Cells.Item(7, 10).Value = 4397.3335
$ws.Cells.Item(7, 11).Value = 3000
$ws.Cells.Item(7, 12).Value = 4397.3335
$ws.Cells.Item(7, 13).Value = -2888
$ws.Cells.Item(7, 14).Value = -4621.3335
$ws.Cells.Item(92, 8).Value = 24999.5
$ws.Cells.Item(92, 10).Value = 24999.5
$ws.Cells.Item(92, 12).Value = 24999.5
$ws.Cells.Item(92, 14).Value = -29991.5
$ws.Cells.Item(126, 8).Value = 4048
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 4397.3335
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 13192.0005
$ws.Cells.Item(126, 13).Value = -6530
$ws.Cells.Item(126, 14).Value = -18132.0005
$ws.Cells.Item(132, 8).Value = 2895.5527
$ws.Cells.Item(132, 9).Value = 2857.0908
$ws.Cells.Item(132, 10).Value = 3149.4
$ws.Cells.Item(132, 11).Value = 8571.2724
$ws.Cells.Item(132, 12).Value = 9448.200000000001
$ws.Cells.Item(132, 13).Value = -6041.2724
$ws.Cells.Item(132, 14).Value = -14508.2
$ws.Cells.Item(136, 8).Value = 1095.4286
$ws.Cells.Item(136, 9).Value = 960.1852
$ws.Cells.Item(136, 10).Value = 1551.875
$ws.Cells.Item(136, 11).Value = 2880.5556
$ws.Cells.Item(136, 12).Value = 4655.625
$ws.Cells.Item(136, 13).Value = -330.5556000000001
$ws.Cells.Item(136, 14).Value = -9755.625
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1944.3636
$ws.Cells.Item(126, 10).Value = 1863.625
$ws.Cells.Item(126, 12).Value = 5590.875
$ws.Cells.Item(126, 14).Value = -10530.875
$ws.Cells.Item(135, 8).Value = 41953
$ws.Cells.Item(135, 10).Value = 41953
$ws.Cells.Item(135, 12).Value = 41953
$ws.Cells.Item(135, 14).Value = -52093
$ws.Cells.Item(136, 8).Value = 658.81134
$ws.Cells.Item(136, 9).Value = 385.375
$ws.Cells.Item(136, 10).Value = 1500.1538
$ws.Cells.Item(136, 11).Value = 1156.125
$ws.Cells.Item(136, 12).Value = 4500.4614
$ws.Cells.Item(136, 13).Value = 1393.875
$ws.Cells.Item(136, 14).Value = -9600.4614

Write-Host "Applied 228 cell updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
